# Update the genetic-algorithm log sheet:
#  - CromosomaMax (col C) changes to the new best chromosome for every data row
#  - Maximo (col D) changes to the new best fitness for every data row
#  - Minimo (col E) / Promedio (col F) change per-row for rows 2-5 (still converging)
#    and then lock onto the new best fitness from row 6 onward (already converged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newChromosome = "110111111111011111011110001011"
$newBest = 0.7654078652278944

# Rows 2-5: unique Minimo/Promedio values while the algorithm is still converging
$minVals = @{
    2 = 0.003768405585372591
    3 = 0.2579061002857124
    4 = 0.2657411921845396
    5 = 0.3904699078424551
}
$avgVals = @{
    2 = 0.2738075225016131
    3 = 0.4548392660494311
    4 = 0.6642810998155401
    5 = 0.7101808155720931
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 201) { $lastRow = 201 }

for ($r = 2; $r -le $lastRow; $r++) {
    # Force text so the 0/1 chromosome string isn't auto-coerced to a number
    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = $newChromosome
    $ws.Cells.Item($r, 4).Value = $newBest

    if ($minVals.ContainsKey($r)) {
        $ws.Cells.Item($r, 5).Value = $minVals[$r]
        $ws.Cells.Item($r, 6).Value = $avgVals[$r]
    } else {
        $ws.Cells.Item($r, 5).Value = $newBest
        $ws.Cells.Item($r, 6).Value = $newBest
    }
}
